# Update countries & provincias Spain
# - Refreshed the "last updated" timestamp
# - Refreshed case/death counters for several countries (rows that moved
#   rank also carry their label along with them)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 00:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 878408
$ws.Range("C4").Value = 29691
$ws.Range("D4").Value = 85549
$ws.Range("E4").Value = 743117
$ws.Range("F4").Value = 14994
$ws.Range("G4").Value = 2083
$ws.Range("H4").Value = 49742

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 153129
$ws.Range("C8").Value = 2481
$ws.Range("D8").Value = 103300
$ws.Range("E8").Value = 44254
$ws.Range("F8").Value = 2908
$ws.Range("G8").Value = 260
$ws.Range("H8").Value = 5575

# --- Row 16: Canada ---
$ws.Range("B16").Value = 42110
$ws.Range("C16").Value = 1920
$ws.Range("D16").Value = 14761
$ws.Range("E16").Value = 25202
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 173
$ws.Range("H16").Value = 2147

# --- Rows 43-44: Chequia / Ucrania swap rank (Chequia moves above Ucrania) ---
$ws.Range("A43").Value = "Chequia"
$ws.Range("B43").Value = 7187
$ws.Range("C43").Value = 55
$ws.Range("D43").Value = 2152
$ws.Range("E43").Value = 4825
$ws.Range("F43").Value = 76
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 210

$ws.Range("A44").Value = "Ucrania"
$ws.Range("B44").Value = 7170
$ws.Range("C44").Value = 578
$ws.Range("D44").Value = 504
$ws.Range("E44").Value = 6479
$ws.Range("F44").Value = 45
$ws.Range("G44").Value = 13
$ws.Range("H44").Value = 187

# --- Row 46: Australia ---
$ws.Range("B46").Value = 6667
$ws.Range("C46").Value = 18
$ws.Range("E46").Value = 1547

# --- Rows 88-89: Nigeria / Tunez swap rank (Nigeria moves above Tunez) ---
$ws.Range("A88").Value = "Nigeria"
$ws.Range("B88").Value = 981
$ws.Range("C88").Value = 108
$ws.Range("D88").Value = 197
$ws.Range("E88").Value = 753
$ws.Range("F88").Value = 2
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 31

$ws.Range("A89").Value = "Tunez"
$ws.Range("B89").Value = 918
$ws.Range("C89").Value = 9
$ws.Range("D89").Value = 190
$ws.Range("E89").Value = 690
$ws.Range("F89").Value = 32
$ws.Range("H89").Value = 38

# --- Row 152 ---
$ws.Range("D152").Value = 30
$ws.Range("E152").Value = 40

# --- Row 155 ---
$ws.Range("B155").Value = 72
$ws.Range("C155").Value = 7
$ws.Range("D155").Value = 14
$ws.Range("E155").Value = 47
$ws.Range("G155").Value = 2
$ws.Range("H155").Value = 11
